# Apply forecast-error-table update:
#  - update existing B2:G10 values (re-run of the underlying forecast-error
#    calculation changed every numeric result, plus the "N" sample-size
#    counts in column G)
#  - add a new row 11 for quarter "Q9" (label in column A, values B11:G11)
#    matching the same formatting as the rows above it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Cells.Item(2, 2).Value = 0.2312766622086449
$ws.Cells.Item(2, 3).Value = 1.879048156589436
$ws.Cells.Item(2, 4).Value = 10.90505104412208
$ws.Cells.Item(2, 5).Value = 3.3022796738196
$ws.Cells.Item(2, 6).Value = 3.326949562682384
$ws.Cells.Item(2, 7).Value = 51

# --- Row 3 ---
$ws.Cells.Item(3, 2).Value = 0.515602234217057
$ws.Cells.Item(3, 3).Value = 1.963007146147895
$ws.Cells.Item(3, 4).Value = 11.82602839579948
$ws.Cells.Item(3, 5).Value = 3.438899300037655
$ws.Cells.Item(3, 6).Value = 3.434545796526614
$ws.Cells.Item(3, 7).Value = 50

# --- Row 4 ---
$ws.Cells.Item(4, 2).Value = 0.3567007860761084
$ws.Cells.Item(4, 3).Value = 1.897578673147804
$ws.Cells.Item(4, 4).Value = 11.03219032630682
$ws.Cells.Item(4, 5).Value = 3.32147411946967
$ws.Cells.Item(4, 6).Value = 3.336486390714962
$ws.Cells.Item(4, 7).Value = 49

# --- Row 5 ---
$ws.Cells.Item(5, 2).Value = 0.5452471532615414
$ws.Cells.Item(5, 3).Value = 1.991328350292171
$ws.Cells.Item(5, 4).Value = 12.46298015738944
$ws.Cells.Item(5, 5).Value = 3.530294627561478
$ws.Cells.Item(5, 6).Value = 3.524844688185481
$ws.Cells.Item(5, 7).Value = 48

# --- Row 6 ---
$ws.Cells.Item(6, 2).Value = 0.4748333535619241
$ws.Cells.Item(6, 3).Value = 1.853041692924432
$ws.Cells.Item(6, 4).Value = 11.1369772540728
$ws.Cells.Item(6, 5).Value = 3.337210999333545
$ws.Cells.Item(6, 6).Value = 3.338969495412706
$ws.Cells.Item(6, 7).Value = 47

# --- Row 7 ---
$ws.Cells.Item(7, 2).Value = 0.4980755295490984
$ws.Cells.Item(7, 3).Value = 1.841574822490943
$ws.Cells.Item(7, 4).Value = 11.41347372914766
$ws.Cells.Item(7, 5).Value = 3.378383301099457
$ws.Cells.Item(7, 6).Value = 3.378389316479261
$ws.Cells.Item(7, 7).Value = 46

# --- Row 8 ---
$ws.Cells.Item(8, 2).Value = 0.3675858156243708
$ws.Cells.Item(8, 3).Value = 1.629918664393685
$ws.Cells.Item(8, 4).Value = 9.8136786659588
$ws.Cells.Item(8, 5).Value = 3.132679151454678
$ws.Cells.Item(8, 6).Value = 3.14619239584994
$ws.Cells.Item(8, 7).Value = 45

# --- Row 9 ---
$ws.Cells.Item(9, 2).Value = 0.3916378424397349
$ws.Cells.Item(9, 3).Value = 1.741856079646775
$ws.Cells.Item(9, 4).Value = 10.36499202449527
$ws.Cells.Item(9, 5).Value = 3.219470767765297
$ws.Cells.Item(9, 6).Value = 3.232505401330167
$ws.Cells.Item(9, 7).Value = 44

# --- Row 10 (F10 previously empty, now populated too) ---
$ws.Cells.Item(10, 2).Value = 0.5067708388377236
$ws.Cells.Item(10, 3).Value = 1.530135071712105
$ws.Cells.Item(10, 4).Value = 9.825491045469901
$ws.Cells.Item(10, 5).Value = 3.134563932267119
$ws.Cells.Item(10, 6).Value = 3.129936092387531
$ws.Cells.Item(10, 7).Value = 43

# --- New row 11 : "Q9" ---
$ws.Cells.Item(11, 1).Value = "Q9"

# Copy the formatting (border/bold/centered style) of the A column label
# cell above onto the new label cell.
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(11, 2).Value = 0.2841523553027562
$ws.Cells.Item(11, 3).Value = 1.690716294712478
$ws.Cells.Item(11, 4).Value = 10.17884259426374
$ws.Cells.Item(11, 5).Value = 3.190429844748782
$ws.Cells.Item(11, 6).Value = 3.216270379803536
$ws.Cells.Item(11, 7).Value = 42
